# Updated cryptos list (price / volume figures refresh + Maker/VeChain row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.745.11'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '3.705.74'
$ws.Range('E3').Value = '  -2.86%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''600.06'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').Value = '''167.80'
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('D7').Value = '3.704.71'
$ws.Range('E7').Value = '  -2.80%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('E10').Value = '  +2.87%  '
$ws.Range('D11').Value = '''6.27'
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').Value = '''0.459'
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('D13').Value = '''38.14'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').Value = '4.324.06'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('D16').Value = '3.701.65'
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').Value = '68.671.77'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '''7.25'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').Value = '''17.20'
$ws.Range('E20').Value = '  +5.59%  '
$ws.Range('D21').Value = '''493.70'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').Value = '''9.19'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('D23').Value = '''0.722'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').Value = '''84.45'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').Value = '''2.30'
$ws.Range('E25').Value = '  -3.88%  '
$ws.Range('E26').Value = '  +1.89%  '
$ws.Range('D27').Value = '''12.21'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').Value = '''10.07'
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '''7.86'
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('D32').Value = '''2.38'
$ws.Range('E32').Value = '  -2.29%  '
$ws.Range('D33').Value = '''31.48'
$ws.Range('E33').Value = '  -3.86%  '
$ws.Range('D34').Value = '3.843.46'
$ws.Range('E34').Value = '  -2.97%  '
$ws.Range('E35').Value = '  -0.58%  '
$ws.Range('D36').Value = '3.646.05'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').Value = '''1.00'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = '''5.75'
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('D40').Value = '''0.133'
$ws.Range('E40').Value = '  -3.19%  '
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('D42').Value = '''49.07'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').Value = '''431.77'
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('D44').Value = '''1.98'
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('D46').Value = '''8.38'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''40.25'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('D49').Value = '''142.18'
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.736.56'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '''0.0349'
$ws.Range('E51').Value = '  -0.56%  '
